$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E3").Value = "test1"
$ws.Range("E4").Value = "test2"

$ws.Range("E5").Select()
